# Presentation outline; fresh analysis results
#
# 1) The auto-updating "datetimeFigureOut" date field (slide master, every
#    slide layout, and the notes master) rolled forward one day:
#    09/03/2023 -> 10/03/2023.
# 2) Slide 6 ("Properties") content placeholder: the trailing colons were
#    dropped from the first two bullets, and two new bullets were added.

$p = $ppt.ActivePresentation

function Update-DatePlaceholder {
    param($container, [string]$newText)

    if ($container -eq $null) { return }
    $shapes = $container.Shapes
    if ($shapes -eq $null) { return }

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        if ($shape.Type -eq 14) {
            if ($shape.PlaceholderFormat.Type -eq 16) {
                $shape.TextFrame.TextRange.Text = $newText
            }
        }
    }
}

$newDate = "10/03/2023"

# Slide master's own date placeholder.
Update-DatePlaceholder $p.SlideMaster $newDate

# Every slide layout hanging off the slide master.
$layouts = $p.SlideMaster.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    Update-DatePlaceholder $layouts.Item($j) $newDate
}

# Notes master's date placeholder too.
if ($p.HasNotesMaster) {
    Update-DatePlaceholder $p.NotesMaster $newDate
}

# Slide 6 ("Properties") - update the bullet list text.
$slide6 = $p.Slides.Item(6)
$contentShape = $slide6.Shapes.Item(2)
$contentShape.TextFrame.TextRange.Text = "Task Formulation`rDisjoint vs Blurry`rOnline vs Offline`rTask-IL vs Class-IL"
